# Generate Report for Handoff
# Adds a new row (row 3) to each of the three worksheets (Overview, zh-cn, de-de)
# describing the hand-off status of a newly generated file:
#   46431e0c-ae49-4950-a8b0-f80923fe2df7oooo....md

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Build the long repeated-'o' filler strings used by the source repo's
# (intentionally long) test file names.
# ---------------------------------------------------------------------------
$o152 = ""
for ($i = 0; $i -lt 152; $i++) { $o152 = $o152 + "o" }

$o37 = ""
for ($i = 0; $i -lt 37; $i++) { $o37 = $o37 + "o" }

$mdFile        = "46431e0c-ae49-4950-a8b0-f80923fe2df7" + $o152 + ".md"
$mdDisplay     = "e2e\" + $mdFile
$xlfZhCn       = "46431e0c-ae49-4950-a8b0-f80923fe2df7" + $o37 + ".696b4cce15046ef8eb7d479cb5622bf94682272d.zh-cn.xlf"
$xlfDeDe       = "46431e0c-ae49-4950-a8b0-f80923fe2df7" + $o37 + ".696b4cce15046ef8eb7d479cb5622bf94682272d.de-de.xlf"

$status        = "Ready for handoff"
$handoffDate   = "2016-10-24 10:39:38"
$handoffDateZh = "2016-10-24 10:39:26"
$handoffDateDe = "2016-10-24 10:39:38"
$neverDate     = "0001-01-01 00:00:00"

$hyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1b06271be099d02e0e02bab908d4ac95a65f4459/e2e/" + $mdFile

# ---------------------------------------------------------------------------
# Sheet "Overview": columns File Name | Path And Name | Extension |
#   Publish URL | zh-cn | de-de | Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $mdFile
$wsOverview.Range("B3").Value = $mdDisplay
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = $status
$wsOverview.Range("F3").Value = $status
$wsOverview.Range("G3").Value = $handoffDate

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $hyperlinkUrl, "", "", $mdDisplay) | Out-Null

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3

# ---------------------------------------------------------------------------
# Sheets "zh-cn" and "de-de": columns Source File Name | File Extension |
#   Status | Source Path | Priority | Content Duplicate | Latest Handoff
#   File | Latest Handoff Datetime | Latest Target File | Latest Handback
#   File | Latest Handback DateTime | Reference Tokens | To be localized |
#   Dependency From | Has metadata | Error Detail
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A3").Value = $mdFile
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = $status
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = $xlfZhCn
$wsZhCn.Range("H3").Value = $handoffDateZh
$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("J3").Value = ""
$wsZhCn.Range("K3").Value = $neverDate
$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "False"
$wsZhCn.Range("P3").Value = ""

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $hyperlinkUrl, "", "", $mdFile) | Out-Null

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P3"))

$wsZhCn.Columns.Item(3).ColumnWidth = 16.3

$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A3").Value = $mdFile
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = $status
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = $xlfDeDe
$wsDeDe.Range("H3").Value = $handoffDateDe
$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("J3").Value = ""
$wsDeDe.Range("K3").Value = $neverDate
$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "False"
$wsDeDe.Range("P3").Value = ""

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $hyperlinkUrl, "", "", $mdFile) | Out-Null

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P3"))

$wsDeDe.Columns.Item(3).ColumnWidth = 16.3

Write-Output "done"
